$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtPaste = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$auto = -4105

# ---------------------------------------------------------------------------
# 1) Materialize the new trailing empty rows (103:119) with the same cell
#    style ("s=4") and row height (23.25) that rows 95:102 already use.
# ---------------------------------------------------------------------------
$ws.Range("A95:D95").Copy()
$ws.Range("A103:D119").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("A103:A119").RowHeight = 23.25

# ---------------------------------------------------------------------------
# 2) Row 95 - cwl_ui_stub_info / CWL 単フレーム情報 / CWL 单帧信息
# ---------------------------------------------------------------------------
$ws.Range("A95").Value = "cwl_ui_stub_info"

$ws.Range("C95").Value = "CWL 単フレーム情報"
$r = $ws.Range("C95").Characters(5, 7)
$r.Font.Name = "宋体"
$r.Font.ColorIndex = $auto

$ws.Range("D95").Value = "CWL 单帧信息"
$r = $ws.Range("D95").Characters(5, 4)
$r.Font.Name = "微软雅黑"
$r.Font.ColorIndex = $auto

# ---------------------------------------------------------------------------
# 3) Row 96 - cwl_ui_stub_header / 占有率 / 平均フレームms / 占比 / 平均帧ms
# ---------------------------------------------------------------------------
$ws.Range("A96").Value = "cwl_ui_stub_header"

$ws.Range("C87").Copy()
$ws.Range("C96").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C96").Value = "占有率 / 平均フレームms"
$r = $ws.Range("C96").Characters(4, 3)
$r.Font.Name = "Cascadia Code"
$r.Font.ColorIndex = $auto
$r = $ws.Range("C96").Characters(7, 6)
$r.Font.Name = "宋体"
$r.Font.ColorIndex = $auto
$r = $ws.Range("C96").Characters(13, 2)
$r.Font.Name = "Cascadia Code"
$r.Font.ColorIndex = $auto

$ws.Range("D65").Copy()
$ws.Range("D96").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D96").Value = "占比 / 平均帧ms"
$r = $ws.Range("D96").Characters(1, 2)
$r.Font.Name = "微软雅黑"
$r.Font.ColorIndex = $auto
$r = $ws.Range("D96").Characters(3, 3)
$r.Font.Name = "Cascadia Code"
$r.Font.ColorIndex = $auto
$r = $ws.Range("D96").Characters(6, 3)
$r.Font.Name = "宋体"
$r.Font.ColorIndex = $auto
$r = $ws.Range("D96").Characters(9, 2)
$r.Font.Name = "Cascadia Code"
$r.Font.ColorIndex = $auto

# ---------------------------------------------------------------------------
# 4) Row 97 - cwl_ui_bgm_view / プレイリストを見る / 查看列表
# ---------------------------------------------------------------------------
$ws.Range("A97").Value = "cwl_ui_bgm_view"

$ws.Range("C87").Copy()
$ws.Range("C97").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C97").Value = "プレイリストを見る"

$ws.Range("D68").Copy()
$ws.Range("D97").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D97").Value = "查看列表"

$ws.Range("A97").EntireRow.RowHeight = 25.5

# ---------------------------------------------------------------------------
# 5) Row 98 - cwl_ui_bgm_hide / プレイリストを隠す / 隐藏列表
# ---------------------------------------------------------------------------
$ws.Range("A98").Value = "cwl_ui_bgm_hide"

$ws.Range("C87").Copy()
$ws.Range("C98").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C98").Value = "プレイリストを隠す"

$ws.Range("D64").Copy()
$ws.Range("D98").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D98").Value = "隐藏列表"

$ws.Range("A98").EntireRow.RowHeight = 25.5

# ---------------------------------------------------------------------------
# 6) Row 99 - cwl_ui_bgm_next / 次の曲 / 下一首
# ---------------------------------------------------------------------------
$ws.Range("A99").Value = "cwl_ui_bgm_next"

$ws.Range("D68").Copy()
$ws.Range("C99").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C99").Value = "次の曲"

$ws.Range("D64").Copy()
$ws.Range("D99").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D99").Value = "下一首"

# ---------------------------------------------------------------------------
# 7) Row 100 - cwl_ui_bgm_last / 前の曲 / 上一首
# ---------------------------------------------------------------------------
$ws.Range("A100").Value = "cwl_ui_bgm_last"

$ws.Range("C87").Copy()
$ws.Range("C100").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C100").Value = "前の曲"

$ws.Range("D64").Copy()
$ws.Range("D100").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D100").Value = "上一首"

$ws.Range("A100").EntireRow.RowHeight = 25.5

# ---------------------------------------------------------------------------
# 8) Row 101 - cwl_ui_bgm_shuffle / シャッフル / 随机播放
# ---------------------------------------------------------------------------
$ws.Range("A101").Value = "cwl_ui_bgm_shuffle"
$ws.Range("C101").Value = "シャッフル"

$ws.Range("D64").Copy()
$ws.Range("D101").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D101").Value = "随机播放"

# ---------------------------------------------------------------------------
# 9) Row 102 - cwl_ui_bgm_rebuild / プレイリストデータの再構築 / 重载列表
# ---------------------------------------------------------------------------
$ws.Range("A102").Value = "cwl_ui_bgm_rebuild"

$ws.Range("C87").Copy()
$ws.Range("C102").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("C102").Value = "プレイリストデータの再構築"

$ws.Range("D68").Copy()
$ws.Range("D102").PasteSpecial($fmtPaste)
$excel.CutCopyMode = $false
$ws.Range("D102").Value = "重载列表"

$ws.Range("A102").EntireRow.RowHeight = 25.5

# ---------------------------------------------------------------------------
# 10) Final view state - matches the author's ending selection/scroll spot.
# ---------------------------------------------------------------------------
$ws.Range("D104").Select()
